# Rubinator3000 PllPatternGenerator - "Make Solving async and work on last layer solving"
#
# The underlying pattern table (columns C/D/E feed formulas in G/H/I/K) gets
# new sample input values in column D, which cascades through the existing
# MOD(...) formulas already in the sheet. Also hide the (now redundant/helper)
# calculation columns F:J so only the human-facing columns show.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Updated "D" (second face) sample inputs for the 4 pattern rows.
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 2

# Hide the helper/lookup columns F:J (formulas in G/H/I and the result in K
# stay live, they're just not meant to be visible to the user anymore).
$ws.Range("F1:J1").EntireColumn.Hidden = $true

$wb.Application.CalculateFull()
